# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: periodo mora changes from 2209 -> 2208 (worker Jose Antonio Perez Diaz)
$ws.Range("E16").Value = "2208"

# Row 17: now holds Ana Maria Ospino Gomez, periodo 2208, updated valores
$ws.Range("C17").Value = "1052989778"
$ws.Range("D17").Value = "ANA MARIA OSPINO GOMEZ"
$ws.Range("F17").Value = 54680
$ws.Range("G17").Value = 1367000

# Row 18: now holds Jose Antonio Perez Diaz, periodo 2209
$ws.Range("C18").Value = "1047366757"
$ws.Range("D18").Value = "JOSE ANTONIO PEREZ DIAZ"
$ws.Range("E18").Value = "2209"
$ws.Range("F18").Value = 62480
$ws.Range("G18").Value = 1562000

# Row 19: stays Ana Maria Ospino Gomez, periodo changes 2210 -> 2209, salario updated
$ws.Range("E19").Value = "2209"
$ws.Range("G19").Value = 1367000

# Row 20: now holds Elena Isabel Garvez Bustillo, periodo 2210, updated valores
$ws.Range("C20").Value = "22534611"
$ws.Range("D20").Value = "ELENA ISABEL GARVEZ BUSTILLO"
$ws.Range("E20").Value = "2210"
$ws.Range("F20").Value = 40000
$ws.Range("G20").Value = 1000000

# Row 21: stays Ana Maria Ospino Gomez, periodo changes 2208 -> 2210, salario updated
$ws.Range("E21").Value = "2210"
$ws.Range("G21").Value = 1367000
